$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "rebs rockfish" sheet: updated GOA-wide harlequin M / tmax AFSC results
# ---------------------------------------------------------------------------
$wsRebs = $wb.Worksheets.Item("rebs rockfish")

$wsRebs.Range("E2").NumberFormat = "@"
$wsRebs.Range("E2").Value = "131"
$wsRebs.Range("F2").NumberFormat = "@"
$wsRebs.Range("F2").Value = "0.0412"

$wsRebs.Range("E3").NumberFormat = "@"
$wsRebs.Range("E3").Value = "120"
$wsRebs.Range("F3").NumberFormat = "@"
$wsRebs.Range("F3").Value = "0.0451"

$wsRebs.Range("E8").NumberFormat = "@"
$wsRebs.Range("E8").Value = "132"
$wsRebs.Range("F8").NumberFormat = "@"
$wsRebs.Range("F8").Value = "0.0409"

$wsRebs.Range("E9").NumberFormat = "@"
$wsRebs.Range("E9").Value = "128"
$wsRebs.Range("F9").NumberFormat = "@"
$wsRebs.Range("F9").Value = "0.0423"

# ---------------------------------------------------------------------------
# 2) "harlequin rockfish" sheet: add two new GOA-wide VBGF rows (10 and 11)
# ---------------------------------------------------------------------------
$wsHarl = $wb.Worksheets.Item("harlequin rockfish")

$wsHarl.Range("A10").Value = "harlequin rockfish"
$wsHarl.Range("B10").Value = "GOA"
$wsHarl.Range("C10").Value = 1
$wsHarl.Range("D10").Value = "VBGF Linf (cm) / k"
$wsHarl.Range("E10").NumberFormat = "@"
$wsHarl.Range("E10").Value = "31.2 / 0.209"
$wsHarl.Range("F10").NumberFormat = "@"
$wsHarl.Range("F10").Value = "0.4212"
$wsHarl.Range("G10").Value = "TenBrink pers comm"

$wsHarl.Range("A11").Value = "harlequin rockfish"
$wsHarl.Range("B11").Value = "GOA"
$wsHarl.Range("C11").Value = 2
$wsHarl.Range("D11").Value = "VBGF Linf (cm) / k"
$wsHarl.Range("E11").NumberFormat = "@"
$wsHarl.Range("E11").Value = "30.9 / 0.167"
$wsHarl.Range("F11").NumberFormat = "@"
$wsHarl.Range("F11").Value = "0.3587"
$wsHarl.Range("G11").Value = "TenBrink pers comm"

# ---------------------------------------------------------------------------
# 3) "rougheye rockfish" sheet: updated tmax AFSC results
# ---------------------------------------------------------------------------
$wsRough = $wb.Worksheets.Item("rougheye rockfish")

$wsRough.Range("E2").NumberFormat = "@"
$wsRough.Range("E2").Value = "116"
$wsRough.Range("F2").NumberFormat = "@"
$wsRough.Range("F2").Value = "0.0466"

$wsRough.Range("E3").NumberFormat = "@"
$wsRough.Range("E3").Value = "96"
$wsRough.Range("F3").NumberFormat = "@"
$wsRough.Range("F3").Value = "0.0565"

$wsRough.Range("E4").NumberFormat = "@"
$wsRough.Range("E4").Value = "107"
$wsRough.Range("F4").NumberFormat = "@"
$wsRough.Range("F4").Value = "0.0505"

$wsRough.Range("E5").NumberFormat = "@"
$wsRough.Range("E5").Value = "104"
$wsRough.Range("F5").NumberFormat = "@"
$wsRough.Range("F5").Value = "0.0518"

$wsRough.Range("E7").NumberFormat = "@"
$wsRough.Range("E7").Value = "113"
$wsRough.Range("F7").NumberFormat = "@"
$wsRough.Range("F7").Value = "0.0476"
